# Applies the "Update state PHEV subsidies (formula fix)" edit to the
# "Passenger Vehicle Calculations" sheet of the BAU Vehicle Subsidy workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Passenger Vehicle Calculations")

# --- Row 47: state subsidy 2000 -> 1000 for years 2020-2024 (D:H) ---
$ws.Range("D47:H47").Value = 1000

# --- Row 50: 2023 (G) becomes a formula tied to row 30, 2024-2026 (H:J) -> 0 ---
$ws.Range("G50").Formula = "=G30*0.5"
$ws.Range("H50:J50").Value = 0

# --- Row 51: 2023-2026 (G:J) 3000 -> 1500 ---
$ws.Range("G51:J51").Value = 1500

# --- Row 52: 2022-2026 (F:J) 1500 -> 1000 ---
$ws.Range("F52:J52").Value = 1000

# --- Row 53: 2023-2026 (G:J) 2000 -> 0 ---
$ws.Range("G53:J53").Value = 0

# --- Row 55: 2020-2026 (D:J) 2200 -> 1000 ---
$ws.Range("D55:J55").Value = 1000

# --- Row 70: fix formula to reference the (now corrected) PHEV subsidy
#     total in row 66 instead of the BEV total in row 65 ---
$ws.Range("C70:AF70").Formula = "=B5+D66"

# --- Restore the saved view state (scroll position / active cell) ---
$ws.Application.ActiveWindow.ScrollRow = $ws.Range("B35").Row
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("B35").Column
$ws.Range("J59").Select()

$wb.Save()
